$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Sheet1: two new rows (10. AHE cl interp buf uchar bank mod)
# ---------------------------------------------------------------------------
$ws1.Range("A22").Value = "10. AHE cl interp buf uchar bank mod"
$ws1.Range("B22").Value = "Histogram calculation"
$ws1.Range("C22").Value = 241.74392

$ws1.Range("A23").Value = "10. AHE cl interp buf uchar bank mod"
$ws1.Range("B23").Value = "Transformation"
$ws1.Range("C23").Value = 163.53752

# row height to match the rest of the sheet
$ws1.Rows.Item(22).RowHeight = 15
$ws1.Rows.Item(23).RowHeight = 15

# C23 uses a distinct (Times 12, black) font/style
$ws1.Range("C23").Font.Name = "Times"
$ws1.Range("C23").Font.Size = 12
$ws1.Range("C23").Font.Color = 0x000000

# ---------------------------------------------------------------------------
# 2. Add Sheet2 right after Sheet1
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Give the whole used range the Sheet1-style font (Times 12, theme text color)
# before filling in values/exceptions, so every cell inherits style "s=1".
$ws2.Range("A1:E11").Font.Name = "Times"
$ws2.Range("A1:E11").Font.Size = 12

# -- Header row --------------------------------------------------------------
$ws2.Range("A1").Value = "Methods"
$ws2.Range("B1").Value = "Hist calculation (ms)"
$ws2.Range("C1").Value = "Transformation (ms)"
$ws2.Range("D1").Value = "Time on original (8MP)"
$ws2.Range("E1").Value = "Time on small (61KP)"

# -- Method names (column A) --------------------------------------------------
$ws2.Range("A2").Value = "1. HE serial"
$ws2.Range("A3").Value = "2. HE numpy"
$ws2.Range("A4").Value = "3. LHE serial"
$ws2.Range("A5").Value = "4. AHE serial"
$ws2.Range("A6").Value = "5. AHE cl buf"
$ws2.Range("A7").Value = "6. AHE cl interp"
$ws2.Range("A8").Value = "7. AHE cl interp buf"
$ws2.Range("A9").Value = "8. AHE cl interp buf uchar"
$ws2.Range("A10").Value = "9. AHE cl interp buf uchar bank"
$ws2.Range("A11").Value = "10. AHE cl interp buf uchar bank mod"

# -- Data values ---------------------------------------------------------
$ws2.Range("D2").Value = 3400
$ws2.Range("E2").Value = 39.4

$ws2.Range("D3").Value = 360
$ws2.Range("E3").Value = 3.25

$ws2.Range("D4").Value = 15000
$ws2.Range("E4").Value = 2790

$ws2.Range("D5").Value = 15000
$ws2.Range("E5").Value = 10500

$ws2.Range("E6").Value = 3.96824

$ws2.Range("B7").Value = 201.055
$ws2.Range("C7").Value = 7.48
$ws2.Range("D7").Value = 208.53
$ws2.Range("E7").Value = 0.49312

$ws2.Range("B8").Value = 203.62
$ws2.Range("C8").Value = 514.47
$ws2.Range("D8").Value = 718.09
$ws2.Range("E8").Value = 3.372

$ws2.Range("B9").Value = 191.29
$ws2.Range("C9").Value = 302.83
$ws2.Range("D9").Value = 494.12
$ws2.Range("E9").Value = 2.542

$ws2.Range("B10").Value = 203.91
$ws2.Range("C10").Value = 289.16
$ws2.Range("D10").Value = 493.08
$ws2.Range("E10").Value = 2.5552

$ws2.Range("B11").Value = 241.74392
$ws2.Range("C11").Value = 399.85936
$ws2.Range("D11").Value = 399.85936
$ws2.Range("E11").Value = 1.94112

# -- Font-color exceptions ----------------------------------------------
# Times 12, black (FF000000) -> style s=3
$blackCells = "D2","D3","D4","E4","D5","E5","D6","E6","D7","C11","E11"
foreach ($addr in $blackCells) {
    $ws2.Range($addr).Font.Name = "Times"
    $ws2.Range($addr).Font.Size = 12
    $ws2.Range($addr).Font.Color = 0x000000
}

# Times 12, dark gray (FF222222) -> style s=4
$ws2.Range("D8").Font.Name = "Times"
$ws2.Range("D8").Font.Size = 12
$ws2.Range("D8").Font.Color = 0x222222

# ---------------------------------------------------------------------------
# 3. Page margins for Sheet2 (match Sheet1: 0.75/0.75/1/1 in, 0.5/0.5 in header/footer)
#    PageSetup margins are expressed in points (72 pt = 1 in).
# ---------------------------------------------------------------------------
$ws2.PageSetup.LeftMargin = 0.75 * 72
$ws2.PageSetup.RightMargin = 0.75 * 72
$ws2.PageSetup.TopMargin = 1 * 72
$ws2.PageSetup.BottomMargin = 1 * 72
$ws2.PageSetup.HeaderMargin = 0.5 * 72
$ws2.PageSetup.FooterMargin = 0.5 * 72

# ---------------------------------------------------------------------------
# 4. Selections: leave Sheet2 selection at I2, then re-select Sheet1 at C23
#    so Sheet1 remains the active/visible tab (matches the authored diff).
# ---------------------------------------------------------------------------
$ws2.Range("I2").Select() | Out-Null
$ws1.Range("C23").Select() | Out-Null
$ws1.Select() | Out-Null
